$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Balance changes to enemy stats (row 2 / enemy 1)
$ws.Range("B2").Value = 24
$ws.Range("D2").Value = 5
$ws.Range("F2").Value = 12

# Balance change to enemy 2's MaxHp
$ws.Range("B3").Value = 48

# Update current selection to match the authored state
$ws.Range("D2").Select()
